# [MOSIP-14369] Fix: boolean values
#
# The "is_active" column (E) used to hold a `=TRUE()` formula that evaluated
# to the number 1. It should instead hold the literal text "TRUE" (a plain
# string, not a boolean/formula). We stage the literal string in a scratch
# cell (built from a formula so the engine treats it as ordinary text instead
# of auto-detecting the keyword TRUE as a boolean the moment it is typed),
# copy it, and paste-special (values only) into each target cell so the
# destination keeps its existing style/number format but picks up a plain
# text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell, far away from the used range, used only to manufacture a
# literal text value "TRUE" (LEFT(...) forces text typing so the engine
# doesn't reinterpret it as a Boolean).
$scratch = $ws.Range("Z100")
$scratch.Formula = '=LEFT("TRUEX",4)'

foreach ($addr in @("E2", "E3", "E4", "E5")) {
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$scratch.ClearContents()
$excel.CutCopyMode = 0

# Update the selection to match: E2:E5 instead of the whole column E.
$ws.Range("E2:E5").Select()
